$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 432.864620784
$ws.Range("D8").Value = 44.05282404768001
$ws.Range("D9").Value = 59.63511917425334
$ws.Range("D10").Value = 207.97596

$ws.Range("D22").Value = 61.91447483243144
$ws.Range("D23").Value = 98.58622320000002
$ws.Range("D24").Value = 72.31800000000001
$ws.Range("D25").Value = 93.22800000000001
$ws.Range("D26").Value = 145.758

$ws.Range("D30").Value = 93.53196000000001
$ws.Range("D31").Value = 32.364721584
